# Commit: "Ajout du mapping 5b9ef178b55ff22959620fdaa372ac3e83db7c97"
# - Bump the Metadata "Date" value to the new generation timestamp.
# - Add a new "Mapping: null" column pair (AK = code, AL = display) to the
#   "Elements" sheet, documenting the DICOM-KOS mapping for each element.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: refresh the generation Date
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-03T11:17:55+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: insert the new Mapping columns AK / AL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Clone the formatting of the last existing column (AJ, "Constraint(s)") onto
# the two new columns so the new cells inherit the same header / body style.
$ws.Range("AJ1:AJ6").Copy()
$ws.Range("AK1:AK6").PasteSpecial(-4122)
$ws.Range("AJ1:AJ6").Copy()
$ws.Range("AL1:AL6").PasteSpecial(-4122)

# -- Header row --------------------------------------------------------------
$ws.Cells.Item(1, 37).Value = "Mapping: null"
$ws.Cells.Item(1, 38).Value = "Mapping: null"

# -- Row 2: Author (root element) --------------------------------------------
$ws.Cells.Item(2, 37).Value = "Author"
$ws.Cells.Item(2, 38).Value = "Author"

# -- Row 3: Author.institution (no DICOM KOS mapping) ------------------------
$ws.Cells.Item(3, 37).Value = ""
$ws.Cells.Item(3, 38).Value = ""

# -- Row 4: Author.person (no DICOM KOS mapping) ------------------------------
$ws.Cells.Item(4, 37).Value = ""
$ws.Cells.Item(4, 38).Value = ""

# -- Row 5: Author.role --------------------------------------------------------
$ws.Cells.Item(5, 37).Value = "author/functionCode@displayName"
$ws.Cells.Item(5, 38).Value = "Cet attribut n'a pas besoin d’être alimenté par un élément du DICOM KOS"

# -- Row 6: Author.specialty ----------------------------------------------------
$ws.Cells.Item(6, 37).Value = "author/assignedAuthor/code@code"
$ws.Cells.Item(6, 38).Value = "Cette métadonnée peut ne pas être renseignée dans le cas d’un DICOM KOS.   Si elle contient une valeur, elle devra contenir le code : 'DISPOSITIF' du JDV_J01_XdsAuthorSpecialty_CISIS"

# -- Column sizing (best-fit widths for the new content) ---------------------
$ws.Columns.Item(37).ColumnWidth = 28.0
$ws.Columns.Item(38).ColumnWidth = 146.16666666666666

Write-Host "Mapping columns added"
